$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tab1")

$data = @(
    @("Thu 06 Feb 2025 02:53:02 PM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Thu 06 Feb 2025 02:56:42 PM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Thu 06 Feb 2025 03:04:47 PM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Thu 06 Feb 2025 03:15:29 PM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Thu 06 Feb 2025 03:17:16 PM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Thu 06 Feb 2025 07:58:59 PM MST`n", "amatos", "Alex", "1", "Lange", "New2"),
    @("Fri 07 Feb 2025 11:43:27 AM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Fri 07 Feb 2025 11:43:42 AM MST`n", "amatos", "Alex", "1", "Lange", "New3"),
    @("Fri 07 Feb 2025 12:08:40 PM MST`n", "amatos", "Alex", "1", "Lange", "New"),
    @("Fri 07 Feb 2025 12:12:28 PM MST`n", "amatos", "Alex", "1", "Lange", "New")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = "'" + $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}
